$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "img"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "rarity"
$ws.Range("D1").Value = "description"
$ws.Range("E1").Value = "location"
$ws.Range("F1").Value = "type"
$ws.Range("G1").Value = "race"
$ws.Range("H1").Value = "bis"

# --- Row 2 (item1 - Pinhole Rifle Widowmaker) ---
$ws.Range("A2").Value = "item1"
$ws.Range("B2").Value = "Pinhole Rifle Widowmaker"
$ws.Range("C2").Value = "purple"
$ws.Range("D2").Value = "trooper,judge"
$ws.Range("E2").Value = "skirmish"
$ws.Range("F2").Value = "nanotechnology"
$ws.Range("G2").Value = "human"
$ws.Range("H2").Value = "'false"

# --- Row 3 (item2 - Biological Transition) ---
$ws.Range("A3").Value = "item2"
$ws.Range("B3").Value = "Biological Transition"
$ws.Range("C3").Value = "purple"
$ws.Range("D3").Value = "assassin,constructor,trooper"
$ws.Range("E3").Value = "skirmish"
$ws.Range("F3").Value = "xenotronics"
$ws.Range("G3").Value = "human"
$ws.Range("H3").Value = "'false"

# --- Row 4 (item2 - Pulson grenade "Doom D3", new row) ---
$ws.Range("A4").Value = "item2"
$ws.Range("B4").Value = "Pulson grenade “Doom D3”"
$ws.Range("C4").Value = "purple"
$ws.Range("D4").Value = "trooper,lord commander"
$ws.Range("E4").Value = "skirmish"
$ws.Range("F4").Value = "xenotronics"
$ws.Range("G4").Value = "human"
$ws.Range("H4").Value = "'false"

$ws.Range("C5").Select()
